$d = $word.ActiveDocument

# Paragraph 1 formatting: add a paragraph border (space-only, no line) and
# bump the left indent from 120 twips (6pt) to 225 twips (11.25pt).
$p1 = $d.Paragraphs(1)
$p1.Format.Borders.DistanceFromTop = 5
$p1.Format.Borders.DistanceFromLeft = 5
$p1.Format.Borders.DistanceFromBottom = 5
$p1.Format.Borders.DistanceFromRight = 5
$p1.Format.LeftIndent = 11.25

# Replace the ID placeholder text (which spans the first run plus a
# trailing-space-only run) with the new ID so the two runs collapse back
# into a single run with no trailing space.
$d.Content.Find.Execute("**ID__AFFARS_pgi_5315_topic_4__ID** ", $false, $false, $false, $false, $false, $true, 1, $false, "**ID__AFFARS_AF_PGI_5315_371_5__ID**", 2)
